# gender_char.xlsx monthly refresh — updates the report-period / download
# timestamp text, plus the refreshed view-count and percentage figures in
# the age/gender breakdown table.
#
# NOTE: every number-looking cell in this sheet (C10:D33) is stored as TEXT
# (inline string), not as a numeric value — that's how the source export
# wrote the workbook. A plain `Range.Value = "86"` assignment would make
# Excel's normal type-inference turn that into a real Number, which would
# not match the original "text" cell type. To avoid that, numeric-looking
# values are entered into a scratch cell with a leading apostrophe (forces
# Excel to keep it as literal text, same as typing '86 into a cell) and
# then copied across with Paste Special (values only) so the destination
# cell's existing style/formatting is left completely untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- plain text fields (not numeric-looking, so no special handling needed) ---
$ws.Range("B4").Value = "2024.10.01. 월간"
$ws.Range("B7").Value = "2024년 11월 10일 24시 38분 00초"

# --- numeric-looking fields that must stay TEXT cells ---
# D8 is a blank cell inside the already-used range, so using it as scratch
# space doesn't grow the sheet's dimension / column list.
$scratch = $ws.Range("D8")

function Set-TextValue($cellRef, $value) {
    $scratch.Value = "'" + $value
    $scratch.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)  # xlPasteValues
}

Set-TextValue "C10" "86"
Set-TextValue "D10" "30.5"

Set-TextValue "C11" "196"
Set-TextValue "D11" "69.5"

Set-TextValue "C14" "3"
Set-TextValue "D14" "1.06"

Set-TextValue "C15" "4"
Set-TextValue "D15" "1.42"

Set-TextValue "C16" "27"
Set-TextValue "D16" "9.57"

Set-TextValue "C17" "32"
Set-TextValue "D17" "11.35"

Set-TextValue "C18" "10"
Set-TextValue "D18" "3.55"

Set-TextValue "C19" "31"
Set-TextValue "D19" "10.99"

Set-TextValue "C20" "4"
Set-TextValue "D20" "1.42"

Set-TextValue "C21" "26"
Set-TextValue "D21" "9.22"

Set-TextValue "C22" "5"
Set-TextValue "D22" "1.77"

Set-TextValue "C23" "16"
Set-TextValue "D23" "5.67"

Set-TextValue "D24" "1.77"

Set-TextValue "C25" "16"
Set-TextValue "D25" "5.67"

Set-TextValue "C26" "4"
Set-TextValue "D26" "1.42"

Set-TextValue "C27" "32"
Set-TextValue "D27" "11.35"

Set-TextValue "C28" "11"
Set-TextValue "D28" "3.9"

Set-TextValue "C29" "24"
Set-TextValue "D29" "8.51"

Set-TextValue "C30" "6"
Set-TextValue "D30" "2.13"

Set-TextValue "C31" "7"
Set-TextValue "D31" "2.48"

Set-TextValue "C32" "11"
Set-TextValue "D32" "3.9"

Set-TextValue "C33" "8"
Set-TextValue "D33" "2.84"

# restore scratch cell to empty/default so it leaves no trace in the sheet
$scratch.Clear()
